$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 24657.928
$ws.Range("I8").Value = 83475.25
$ws.Range("K8").Value = 250425.75
$ws.Range("M8").Value = -250286.75

$ws.Range("H21").Value = 39969.5
$ws.Range("I21").Value = 39969.5
$ws.Range("K21").Value = 39969.5
$ws.Range("M21").Value = -39501.5

$ws.Range("H23").Value = 39969.5
$ws.Range("I23").Value = 39969.5
$ws.Range("K23").Value = 39969.5
$ws.Range("M23").Value = -39735.5

$ws.Range("H38").Value = 1601
$ws.Range("I38").Value = 619.4545000000001
$ws.Range("K38").Value = 1858.3635
$ws.Range("M38").Value = -1486.3635

$ws.Range("H76").Value = 91014200
$ws.Range("I76").Value = 134648.75
$ws.Range("K76").Value = 134648.75
$ws.Range("M76").Value = -134333.75

$ws.Range("H79").Value = 91014200
$ws.Range("I79").Value = 134648.75
$ws.Range("K79").Value = 134648.75
$ws.Range("M79").Value = -133556.75

$ws.Range("H112").Value = 2425.205
$ws.Range("J112").Value = 2463.2632
$ws.Range("L112").Value = 7389.7896
$ws.Range("N112").Value = -9605.7896

$ws.Range("H113").Value = 4643.8
$ws.Range("I113").Value = 4999
$ws.Range("K113").Value = 4999
$ws.Range("M113").Value = -1745

$ws.Range("H138").Value = 6555.8223
$ws.Range("I138").Value = 5782.3335
$ws.Range("K138").Value = 17347.0005
$ws.Range("M138").Value = -12207.0005

$ws.Range("H141").Value = 1796.3478
$ws.Range("J141").Value = 6125
$ws.Range("L141").Value = 18375
$ws.Range("N141").Value = -28735

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1256
$ws.Range("I2").Value = 1299.5
$ws.Range("K2").Value = 1299.5
$ws.Range("M2").Value = -1186.5

$ws.Range("H116").Value = 1256
$ws.Range("I116").Value = 1299.5
$ws.Range("K116").Value = 1299.5
$ws.Range("M116").Value = 994.5

$ws.Range("H132").Value = 5619.6943
$ws.Range("I132").Value = 1652.909
$ws.Range("K132").Value = 4958.727000000001
$ws.Range("M132").Value = -2428.727000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1256
$ws.Range("I3").Value = 1299.5
$ws.Range("K3").Value = 1299.5
$ws.Range("M3").Value = -1185.5

$ws.Range("H20").Value = 2826.647
$ws.Range("J20").Value = 3631.875
$ws.Range("L20").Value = 3631.875
$ws.Range("N20").Value = -4125.875

$ws.Range("H94").Value = 1047.2
$ws.Range("I94").Value = 941.3333
$ws.Range("K94").Value = 941.3333
$ws.Range("M94").Value = -490.3333

$ws.Range("H134").Value = 6703.9375
$ws.Range("I134").Value = 10100
$ws.Range("J134").Value = 6218.7856
$ws.Range("K134").Value = 30300
$ws.Range("L134").Value = 18656.3568
$ws.Range("M134").Value = -27765
$ws.Range("N134").Value = -23726.3568

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 656.2353000000001
$ws.Range("I22").Value = 515.36365
$ws.Range("J22").Value = 914.5
$ws.Range("K22").Value = 515.36365
$ws.Range("L22").Value = 914.5
$ws.Range("M22").Value = -165.36365
$ws.Range("N22").Value = -1614.5

$ws.Range("H109").Value = 74986
$ws.Range("J109").Value = 74986
$ws.Range("L109").Value = 74986
$ws.Range("N109").Value = -77066

$ws.Range("H132").Value = 3134.2092
$ws.Range("I132").Value = 1867.9656
$ws.Range("J132").Value = 5757.143
$ws.Range("K132").Value = 5603.8968
$ws.Range("L132").Value = 17271.429
$ws.Range("M132").Value = -3073.8968
$ws.Range("N132").Value = -22331.429

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 11000001
$ws.Range("I32").Value = 30000000
$ws.Range("K32").Value = 90000000
$ws.Range("M32").Value = -89999717

$ws.Range("H34").Value = 72374
$ws.Range("J34").Value = 144687.5
$ws.Range("L34").Value = 434062.5
$ws.Range("N34").Value = -434230.5

$ws.Range("H39").Value = 15404.223
$ws.Range("J39").Value = 19520
$ws.Range("L39").Value = 58560
$ws.Range("N39").Value = -59148

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 15043.8
$ws.Range("I43").Value = 5828.5713
$ws.Range("J43").Value = 23107.125
$ws.Range("K43").Value = 5828.5713
$ws.Range("L43").Value = 23107.125
$ws.Range("M43").Value = -5677.5713
$ws.Range("N43").Value = -23409.125

$ws.Range("H70").Value = 50007084
$ws.Range("I70").Value = 5952.75
$ws.Range("K70").Value = 5952.75
$ws.Range("M70").Value = -5682.75

$ws.Range("H73").Value = 50007084
$ws.Range("I73").Value = 5952.75
$ws.Range("K73").Value = 5952.75
$ws.Range("M73").Value = -5016.75

$ws.Range("H80").Value = 2226007.2
$ws.Range("I80").Value = 1254257.4
$ws.Range("J80").Value = 10000006
$ws.Range("K80").Value = 1254257.4
$ws.Range("L80").Value = 10000006
$ws.Range("M80").Value = -1253259.4
$ws.Range("N80").Value = -10002002

$ws.Range("H83").Value = 2226007.2
$ws.Range("I83").Value = 1254257.4
$ws.Range("J83").Value = 10000006
$ws.Range("K83").Value = 6271287
$ws.Range("L83").Value = 50000030
$ws.Range("M83").Value = -6266295
$ws.Range("N83").Value = -50010014

$ws.Range("H132").Value = 274242.2
$ws.Range("I132").Value = 337053.97
$ws.Range("K132").Value = 1011161.91
$ws.Range("M132").Value = -1008631.91

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2503788
$ws.Range("I40").Value = 3336251
$ws.Range("K40").Value = 3336251
$ws.Range("M40").Value = -3336115

$ws.Range("H132").Value = 4916.5557
$ws.Range("I132").Value = 4071.4285
$ws.Range("K132").Value = 12214.2855
$ws.Range("M132").Value = -9684.2855

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 12457.546
$ws.Range("J81").Value = 31983.75
$ws.Range("L81").Value = 63967.5
$ws.Range("N81").Value = -66089.5

$ws.Range("H84").Value = 12457.546
$ws.Range("J84").Value = 31983.75
$ws.Range("L84").Value = 319837.5
$ws.Range("N84").Value = -330445.5

$ws.Range("H122").Value = 31254504
$ws.Range("I122").Value = 50003360
$ws.Range("J122").Value = 6416.1665
$ws.Range("K122").Value = 150010080
$ws.Range("L122").Value = 19248.4995
$ws.Range("M122").Value = -150007630
$ws.Range("N122").Value = -24148.4995

$ws.Range("H132").Value = 3583.8157
$ws.Range("I132").Value = 2493.7778
$ws.Range("J132").Value = 6259.364
$ws.Range("K132").Value = 7481.3334
$ws.Range("L132").Value = 18778.092
$ws.Range("M132").Value = -4951.3334
$ws.Range("N132").Value = -23838.092
